{"js": "// Apply the benchmark-table updates described by the diff.\n// The document contains a single table, one value per row (one cell per row).\n// Changes:\n//  - Several simple text replacements (row index -> old/new value).\n//  - Three rows whose content was a tab-separated list of numbers; their\n//    entire content collapses into a single value (the value previously\n//    held by one of the simple rows near the top of the table).\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// 0-based row index -> new text for that row's single cell.\nconst updates = {\n  0: \"0M\",        // was 94.54\n  1: \"0M\",        // was 1.4\n  2: \"0M\",        // was 25\n  3: \"117\",       // was 84\n  5: \"0.34491\",   // was 0.00050\n  6: \"0.05456\",   // was 0.00013\n  7: \"0.00849\",   // was 0.00005\n  8: \"0.33389\",   // was 0.00014\n  9: \"0.33389\",   // was 0.00017\n  10: \"0.34491\",  // was 0.00022\n  11: \"0.78895\",  // was 0.01248\n  43: \"94.54\",    // was \"2\\t0.33389\\t0.34491\\t0.33940\\t0.00780\\t0.33389\\t0.33389\\t0.34491\\t0.67880\\t56.4\"\n  44: \"1.4\",      // was \"2\\t0.00529\\t0.07772\\t0.04151\\t0.05122\\t0.00529\\t0.00529\\t0.07772\\t0.08301\\t56.4\"\n  45: \"25\",       // was \"29\\t0.00017\\t0.00077\\t0.00051\\t0.00018\\t0.00035\\t0.00054\\t0.00065\\t0.01466\\t56.4\"\n};\n\nfor (const rowIndexStr of Object.keys(updates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const newText = updates[rowIndex];\n  const cell = table.getCell(rowIndex, 0);\n  const range = cell.body.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-table updates described by the diff.\n# The document contains a single table, one value per row (one cell per row).\n# Changes:\n#  - Several simple text replacements (row index -> old/new value).\n#  - Three rows whose content was a tab-separated list of numbers; their\n#    entire content collapses into a single value (the value previously\n#    held by one of the simple rows near the top of the table).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# 1-based row number -> new text for that row's single cell (column 1).\n$updates = [ordered]@{\n    1  = \"0M\"       # was 94.54\n    2  = \"0M\"       # was 1.4\n    3  = \"0M\"       # was 25\n    4  = \"117\"      # was 84\n    6  = \"0.34491\"  # was 0.00050\n    7  = \"0.05456\"  # was 0.00013\n    8  = \"0.00849\"  # was 0.00005\n    9  = \"0.33389\"  # was 0.00014\n    10 = \"0.33389\"  # was 0.00017\n    11 = \"0.34491\"  # was 0.00022\n    12 = \"0.78895\"  # was 0.01248\n    44 = \"94.54\"    # was \"2`t0.33389`t0.34491`t0.33940`t0.00780`t0.33389`t0.33389`t0.34491`t0.67880`t56.4\"\n    45 = \"1.4\"      # was \"2`t0.00529`t0.07772`t0.04151`t0.05122`t0.00529`t0.00529`t0.07772`t0.08301`t56.4\"\n    46 = \"25\"       # was \"29`t0.00017`t0.00077`t0.00051`t0.00018`t0.00035`t0.00054`t0.00065`t0.01466`t56.4\"\n}\n\nforeach ($rowNum in $updates.Keys) {\n    $cell = $tbl.Cell($rowNum, 1)\n    $cell.Range.Text = $updates[$rowNum]\n}\n"}
